# Updated cryptos list values (Price / Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.206.28"
$ws.Range("E2").Value = "  +5.76%  "
$ws.Range("D3").Value = "1.782.84"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'243.97"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4916"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.2668"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.06249"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "1.775.92"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").Value = "'16.48"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "'0.07013"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'0.6263"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").Value = "'4.634"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "28.171.44"
$ws.Range("E16").Value = "  +6.34%  "
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'0.9996"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'0.000007227"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'12.06"
$ws.Range("E20").Value = "  +5.99%  "
$ws.Range("D21").Value = "2.005.48"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'4.560"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").Value = "'8.733"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").Value = "'5.232"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").Value = "'141.25"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("D27").Value = "'1.858"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").Value = "'109.10"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").Value = "'1.385"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "'4.174"
$ws.Range("E30").Value = "  +6.70%  "
$ws.Range("D31").Value = "'0.08242"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "'3.767"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "'0.04885"
$ws.Range("E33").Value = "  +8.95%  "
$ws.Range("D34").Value = "'1.072"
$ws.Range("E34").Value = "  +7.11%  "
$ws.Range("D35").Value = "'2.612"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +4.30%  "
$ws.Range("D37").Value = "'0.9463"
$ws.Range("D38").Value = "'2.589"
$ws.Range("E38").Value = "  +7.46%  "
$ws.Range("D39").Value = "'2.051"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'5.923"
$ws.Range("E40").Value = "  +5.80%  "
$ws.Range("D41").Value = "'0.01550"
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "'99.44"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.3975"
$ws.Range("E44").Value = "  +3.27%  "
$ws.Range("D45").Value = "'7.182"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").Value = "'0.1203"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("D47").Value = "'0.05427"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "'7.995"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("D50").Value = "'30.61"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "'52.82"
$ws.Range("E51").Value = "  +2.49%  "
